$wb = $excel.ActiveWorkbook

$wsNote = $wb.Worksheets.Item("Note")
$wsDevice = $wb.Worksheets.Item("device")

# The free-format note rows (title/update-date/comment) at the top of the
# "device" sheet duplicated what already lives on the "Note" sheet, and
# they broke the "all sheets -> one YAML" export, so drop them. Everything
# below shifts up (header + the two device rows now start at row 1).
$wsDevice.Rows("1:3").Delete()

# Update the remembered selection on the "Note" sheet.
$wsNote.Activate()
$wsNote.Range("C33").Select()

# "device" becomes the active tab / selected sheet of the workbook.
$wsDevice.Activate()
$wsDevice.Range("A1").Select()
